$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("all_systems")

$ws.Range("AB3").Value = 37
$ws.Range("AB4").Value = 38
$ws.Range("AB5").Value = 39
$ws.Range("AB6").Value = 40
$ws.Range("AB7").Value = 41
$ws.Range("AB8").Value = 42
$ws.Range("AB9").Value = 43
$ws.Range("AB10").Value = 44
$ws.Range("AB12").Value = 45
$ws.Range("AB13").Value = 46
$ws.Range("AB34").Value = 26
$ws.Range("AB35").Value = 27
$ws.Range("AB36").Value = 28
$ws.Range("AB37").Value = 29
$ws.Range("AB38").Value = 30
$ws.Range("AB39").Value = 31
$ws.Range("AB40").Value = 32
$ws.Range("AB41").Value = 33
$ws.Range("AB42").Value = 34
$ws.Range("AB43").Value = 35
$ws.Range("AB44").Value = 36

$ws.Activate()
$ws.Application.ActiveWindow.SplitRow = 2
$ws.Application.ActiveWindow.SplitColumn = 1
$ws.Range("R3").Select()
$ws.Application.ActiveWindow.FreezePanes = $true
$ws.Range("AB14").Select()
